$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gchart")

$ws.Range("B60").Value = "28: [Asc] "
$ws.Range("B67").Value = "19: [Mars] "
$ws.Range("B68").Value = "21: [Saturn] "
$ws.Range("B69").Value = "26: [Pluto] "
$ws.Range("B74").Value = "10: [Jupiter] "
$ws.Range("B82").Value = "4: [Uranus] "
$ws.Range("B89").Value = "26: [Neptune] "
$ws.Range("B96").Value = "17: [Moon] "
$ws.Range("B99").Value = "25: [Venus] "
$ws.Range("B102").Value = "5: [Mercury] "
$ws.Range("B111").Value = "1: [Sun] "
